$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - best_params (inline strings)
$ws.Range("B2").Value = "{'alpha': 0.01, 'max_iter': 1000}"
$ws.Range("E2").Value = "{'C': 100, 'gamma': 1}"
$ws.Range("F2").Value = "{'max_depth': 20, 'n_estimators': 10}"
$ws.Range("G2").Value = "{'learning_rate': 0.1, 'max_depth': 3, 'n_estimators': 200}"
$ws.Range("H2").Value = "{'learning_rate': 1, 'n_estimators': 10}"
$ws.Range("K2").Value = "{'activation': 'leaky_relu', 'b_random_vec_range': [0, 10], 'lam': 1, 'n_layer': 32, 'n_nodes': 256, 'random_seed': 542, 'same_feature': True, 'w_random_vec_range': [-10, 10]}"

# Row 3 - rmse
$ws.Range("B3").Value = 0.07201393407362781
$ws.Range("C3").Value = 0.09460945310694435
$ws.Range("D3").Value = 0.08259024985382742
$ws.Range("E3").Value = 0.07609899163810305
$ws.Range("F3").Value = 0.05884017133411264
$ws.Range("G3").Value = 0.05563437495359801
$ws.Range("H3").Value = 0.10547618456903
$ws.Range("I3").Value = 0.04820478403272017
$ws.Range("J3").Value = 0.05565852779068281
$ws.Range("K3").Value = 0.02602616318137193

# Row 4 - r2
$ws.Range("B4").Value = 0.9102808336638673
$ws.Range("C4").Value = 0.8454503780815263
$ws.Range("D4").Value = 0.882703836767836
$ws.Range("E4").Value = 0.8998901436423736
$ws.Range("F4").Value = 0.9400765369332685
$ws.Range("G4").Value = 0.946519544676196
$ws.Range("H4").Value = 0.8093129765593445
$ws.Range("I4").Value = 0.9588445467502249
$ws.Range("J4").Value = 0.9457495075857937
$ws.Range("K4").Value = 0.9875439657656241

# Row 5 - mape
$ws.Range("B5").Value = 12.78659450029113
$ws.Range("C5").Value = 17.80408635620394
$ws.Range("D5").Value = 15.22375589511617
$ws.Range("E5").Value = 14.9100138704676
$ws.Range("F5").Value = 7.419014189760508
$ws.Range("G5").Value = 7.881477329098279
$ws.Range("H5").Value = 22.71974804164502
$ws.Range("I5").Value = 6.645592447108496
$ws.Range("J5").Value = 7.748171469649608
$ws.Range("K5").Value = 3.964698070926814

$wb.Save()
